$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("File")
Write-Host $ws.Name
